# Updated queries for C3DC first half testcases.
# This script rewrites the embedded SQL queries on Sheet1 so that the
# join conditions use the renamed id columns (study_id / participant_id)
# instead of the old generic "id" columns, and normalizes a couple of
# "WHERE " lines that lost their trailing space.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Fix-Query([string]$text, [bool]$trimWhere) {
    $t = $text
    $t = $t.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $t = $t.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $t = $t.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $t = $t.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $t = $t.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $t = $t.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    if ($trimWhere) {
        $t = $t.Replace("WHERE `r`n", "WHERE`r`n")
        $t = $t.Replace("WHERE `n", "WHERE`n")
    }
    return $t
}

# Cell C2: COUNT/summary query - "WHERE " stays as-is (trailing space kept)
$ws.Range("C2").Value = (Fix-Query $ws.Range("C2").Value() $false)

# Cell B2: dbGaP accession / study name query - trailing space on WHERE removed
$ws.Range("B2").Value = (Fix-Query $ws.Range("B2").Value() $true)

# Cell B3: Participant Id query - trailing space on WHERE removed
$ws.Range("B3").Value = (Fix-Query $ws.Range("B3").Value() $true)

# Cell B4: Diagnosis query - trailing space on WHERE removed
$ws.Range("B4").Value = (Fix-Query $ws.Range("B4").Value() $true)

# Cell B5: Treatment query - trailing space on WHERE removed
$ws.Range("B5").Value = (Fix-Query $ws.Range("B5").Value() $true)

# Cell B6: Treatment Response query - trailing space on WHERE removed
$ws.Range("B6").Value = (Fix-Query $ws.Range("B6").Value() $true)

# Cell B7: Survival query - "WHERE " stays as-is (trailing space kept)
$ws.Range("B7").Value = (Fix-Query $ws.Range("B7").Value() $false)

# Column C width grew (and lost its "best fit" flag) now that the query text changed.
# Note: Excel's ColumnWidth property (character units) and the raw OOXML "width"
# attribute differ by a small fixed padding offset (~5/6 of a character here), so
# we back-solve the ColumnWidth value that serializes to width="68.5" in the XML.
$ws.Columns.Item(3).ColumnWidth = 67.6666666666667
